# "13 - Subprograms" deck: minor wording fixes on the "Converting NamedValue
# to Variable (continued)" slide (slide 43).
#
#   1. "    ...  // check that named values are being passed for"
#        -> "    ...  // check that named values are being passed"
#   2. "         //    var parameters (see next slide)"
#        -> "         // for var parameters (see next slide)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(43)
$shp = $s.Shapes.Item(2)   # "Content Placeholder 2"
$tr = $shp.TextFrame.TextRange

# --- Fix 1: drop the trailing " for" from the "check that named values..." line.
$old1 = "    ...  // check that named values are being passed for"
$new1 = "    ...  // check that named values are being passed"
$text = $tr.Text
$idx1 = $text.IndexOf($old1)
$run1 = $tr.Characters($idx1 + 1, $old1.Length)
$run1.Text = $new1

# --- Fix 2: merge the "//    " / "var" / " parameters (see next slide)" runs
#     into a single run that reads "// for var parameters (see next slide)".
$old2 = "         //    var parameters (see next slide)"
$new2 = "         // for var parameters (see next slide)"
$text = $tr.Text
$idx2 = $text.IndexOf($old2)
$run2 = $tr.Characters($idx2 + 1, $old2.Length)
$run2.Text = $new2
